$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1904761904761905
$ws.Range("C2").Value = 0.5854341736694678
$ws.Range("J2").Value = 0.01680672268907563
$ws.Range("P2").Value = 0.134453781512605
$ws.Range("S2").Value = 0.07282913165266107
$ws.Range("C3").Value = 0.04504504504504504
$ws.Range("J3").Value = 0.004504504504504504
$ws.Range("P3").Value = 0.7207207207207207
$ws.Range("S3").Value = 0.2297297297297297
$ws.Range("J4").Value = 0.03571428571428571
$ws.Range("P4").Value = 0.6607142857142857
$ws.Range("S4").Value = 0.3035714285714285
$ws.Range("B6").Value = 0.0532319391634981
$ws.Range("D6").Value = 0.01140684410646388
$ws.Range("F6").Value = 0.05703422053231939
$ws.Range("J6").Value = 0.2395437262357414
$ws.Range("O6").Value = 0.03041825095057034
$ws.Range("Q6").Value = 0.1939163498098859
$ws.Range("R6").Value = 0.04562737642585551
$ws.Range("S6").Value = 0.3688212927756654
$ws.Range("B7").Value = 0.08230452674897119
$ws.Range("D7").Value = 0.01646090534979424
$ws.Range("F7").Value = 0.05349794238683128
$ws.Range("J7").Value = 0.1440329218106996
$ws.Range("O7").Value = 0.02880658436213992
$ws.Range("Q7").Value = 0.1563786008230453
$ws.Range("R7").Value = 0.05761316872427984
$ws.Range("S7").Value = 0.4609053497942387
$ws.Range("B8").Value = 0.1315240083507307
$ws.Range("D8").Value = 0.01461377870563674
$ws.Range("E8").Value = 0.00208768267223382
$ws.Range("F8").Value = 0.04175365344467641
$ws.Range("J8").Value = 0.09812108559498957
$ws.Range("O8").Value = 0.01878914405010438
$ws.Range("Q8").Value = 0.1920668058455115
$ws.Range("R8").Value = 0.06680584551148225
$ws.Range("S8").Value = 0.4342379958246347
$ws.Range("B9").Value = 0.1092436974789916
$ws.Range("D9").Value = 0.02100840336134454
$ws.Range("F9").Value = 0.06302521008403361
$ws.Range("J9").Value = 0.1260504201680672
$ws.Range("O9").Value = 0.02941176470588235
$ws.Range("Q9").Value = 0.1470588235294118
$ws.Range("R9").Value = 0.05882352941176471
$ws.Range("S9").Value = 0.4453781512605042
$ws.Range("B10").Value = 0.1231231231231231
$ws.Range("D10").Value = 0.02927927927927928
$ws.Range("E10").Value = 0.0007507507507507507
$ws.Range("F10").Value = 0.08483483483483484
$ws.Range("J10").Value = 0.1013513513513514
$ws.Range("O10").Value = 0.02102102102102102
$ws.Range("Q10").Value = 0.1951951951951952
$ws.Range("R10").Value = 0.0518018018018018
$ws.Range("S10").Value = 0.3926426426426426
$ws.Range("G11").Value = 0.1559139784946237
$ws.Range("J11").Value = 0.08333333333333333
$ws.Range("K11").Value = 0.2016129032258064
$ws.Range("L11").Value = 0.5456989247311828
$ws.Range("S11").Value = 0.01344086021505376
$ws.Range("G12").Value = 0.7451923076923077
$ws.Range("J12").Value = 0.2067307692307692
$ws.Range("K12").Value = 0.004807692307692308
$ws.Range("L12").Value = 0.01442307692307692
$ws.Range("S12").Value = 0.02884615384615385
$ws.Range("G13").Value = 0.7659574468085106
$ws.Range("J13").Value = 0.2340425531914894
$ws.Range("F15").Value = 0.02489626556016597
$ws.Range("H15").Value = 0.1161825726141079
$ws.Range("I15").Value = 0.08298755186721991
$ws.Range("J15").Value = 0.3236514522821577
$ws.Range("K15").Value = 0.09958506224066389
$ws.Range("M15").Value = 0.01244813278008299
$ws.Range("O15").Value = 0.04149377593360996
$ws.Range("S15").Value = 0.2987551867219917
$ws.Range("F16").Value = 0.01255230125523013
$ws.Range("H16").Value = 0.1422594142259414
$ws.Range("I16").Value = 0.1129707112970711
$ws.Range("J16").Value = 0.401673640167364
$ws.Range("K16").Value = 0.1171548117154812
$ws.Range("M16").Value = 0.01255230125523013
$ws.Range("O16").Value = 0.04602510460251046
$ws.Range("S16").Value = 0.1548117154811715
$ws.Range("F17").Value = 0.02320675105485232
$ws.Range("H17").Value = 0.1877637130801688
$ws.Range("I17").Value = 0.1075949367088608
$ws.Range("J17").Value = 0.4177215189873418
$ws.Range("K17").Value = 0.08649789029535865
$ws.Range("M17").Value = 0.0189873417721519
$ws.Range("O17").Value = 0.04219409282700422
$ws.Range("S17").Value = 0.1160337552742616
$ws.Range("F18").Value = 0.0352112676056338
$ws.Range("H18").Value = 0.2183098591549296
$ws.Range("I18").Value = 0.1056338028169014
$ws.Range("J18").Value = 0.3661971830985916
$ws.Range("K18").Value = 0.1126760563380282
$ws.Range("M18").Value = 0.02112676056338028
$ws.Range("O18").Value = 0.06338028169014084
$ws.Range("S18").Value = 0.07746478873239436
$ws.Range("F19").Value = 0.02059025394646534
$ws.Range("H19").Value = 0.209334248455731
$ws.Range("I19").Value = 0.08647906657515443
$ws.Range("J19").Value = 0.3603294440631434
$ws.Range("K19").Value = 0.1276595744680851
$ws.Range("M19").Value = 0.02127659574468085
$ws.Range("O19").Value = 0.07343857240905971
$ws.Range("S19").Value = 0.1008922443376802
